$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.368.02"
$ws.Range("E2").Value = "  +2.32%  "

$ws.Range("D3").Value = "2.593.67"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'531.31"
$ws.Range("E5").Value = "  +2.71%  "

$ws.Range("D6").Value = "'140.74"
$ws.Range("E6").Value = "  +1.19%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("D9").Value = "2.606.64"
$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("D10").Value = "'6.46"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  +2.42%  "

$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "3.051.42"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "59.271.36"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.647.84"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000134"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "'346.46"

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").Value = "'10.10"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "'6.38"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'67.62"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").Value = "'0.405"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "'7.15"
$ws.Range("E28").Value = "  +3.12%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").Value = "0.0₃0733"
$ws.Range("E30").Value = "  +2.23%  "

$ws.Range("E31").Value = "  +3.85%  "

$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("D33").Value = "'18.77"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").Value = "'149.65"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").Value = "'3.96"
$ws.Range("E35").Value = "  +1.09%  "

$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("D37").Value = "'36.85"
$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("E38").Value = "  +3.26%  "

$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "'0.832"
$ws.Range("E39").Value = "  +0.98%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'0.835"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").Value = "'271.63"
$ws.Range("E43").Value = "  -0.79%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'10.78"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.596"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  +1.89%  "

$ws.Range("D47").Value = "'0.0519"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.945.93"
$ws.Range("E48").Value = "  -1.45%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0221"
$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'18.26"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.47"
$ws.Range("E51").Value = "  -0.66%  "
